# "Generate Report for Handoff"
#
# The localization-status report is regenerated: the "Handed back: in sync
# with en-US" status becomes "Ready for handoff" everywhere it appears
# (Overview!E2, Overview!F2, zh-cn!C2, de-de!C2), and the associated
# "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps are
# refreshed to the new handoff run's timestamps. The Status/zh-cn/de-de
# columns also get narrower now that the new status text is shorter than
# the old one.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff" ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value     = "Ready for handoff"
$wsDeDe.Range("C2").Value     = "Ready for handoff"

# --- Refreshed timestamps for the new handoff ---
$wsOverview.Range("G2").Value = "2016-08-20 23:03:30"
$wsZhCn.Range("H2").Value     = "2016-08-20 23:03:26"
$wsDeDe.Range("H2").Value     = "2016-08-20 23:03:30"

# --- Column widths shrink to fit the new, shorter status text ---
# (29.9777047293527 characters -> 17.2159881591797 characters)
$newStatusColWidth = 16.333333333333332
$wsOverview.Columns.Item(5).ColumnWidth = $newStatusColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newStatusColWidth
$wsZhCn.Columns.Item(3).ColumnWidth     = $newStatusColWidth
$wsDeDe.Columns.Item(3).ColumnWidth     = $newStatusColWidth
